$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.384.70'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '3.318.90'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'585.54"
$ws.Range('E5').Value = '  +2.02%  '
$ws.Range('D6').Value = "'180.77"
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').Value = "'0.653"
$ws.Range('E7').Value = '  +5.89%  '
$ws.Range('D9').Value = '3.318.51'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').Value = "'0.127"
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('D11').Value = "'6.82"
$ws.Range('E11').Value = '  +2.61%  '
$ws.Range('D12').Value = "'0.402"
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = '3.897.77'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('E14').Value = '  -2.86%  '
$ws.Range('D15').Value = '66.406.48'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = "'0.0000164"
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.318.49'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = "'424.61"
$ws.Range('E19').Value = '  -2.73%  '
$ws.Range('E20').Value = '  -3.05%  '
$ws.Range('D21').Value = "'13.16"
$ws.Range('E21').Value = '  -2.89%  '
$ws.Range('D22').Value = "'7.38"
$ws.Range('E22').Value = '  -2.60%  '
$ws.Range('D23').Value = "'71.83"
$ws.Range('E23').Value = '  -1.86%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = "'5.68"
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').Value = '3.472.53'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').Value = "'0.514"
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('E28').Value = '  +6.90%  '
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('D30').Value = "'9.05"
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('D33').Value = "'22.39"
$ws.Range('E33').Value = '  -1.48%  '
$ws.Range('D35').Value = "'5.18"
$ws.Range('E35').Value = '  -1.29%  '
$ws.Range('D36').Value = "'6.63"
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('D38').Value = "'160.25"
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').Value = "'1.44"
$ws.Range('E39').Value = '  -2.76%  '
$ws.Range('D40').Value = '2.865.39'
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('D41').Value = "'1.81"
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('D42').Value = "'26.41"
$ws.Range('E42').Value = '  -4.93%  '
$ws.Range('D43').Value = "'4.33"
$ws.Range('E43').Value = '  -2.30%  '
$ws.Range('D44').Value = "'0.760"
$ws.Range('E44').Value = '  -4.02%  '
$ws.Range('D45').Value = "'39.77"
$ws.Range('E45').Value = '  -1.08%  '
$ws.Range('D46').Value = "'0.0662"
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('D47').Value = "'5.93"
$ws.Range('E47').Value = '  -4.16%  '
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('D49').Value = "'23.15"
$ws.Range('E49').Value = '  -3.91%  '
$ws.Range('D50').Value = "'312.11"
$ws.Range('E50').Value = '  -3.52%  '
$ws.Range('D51').Value = "'0.0273"
$ws.Range('E51').Value = '  +0.84%  '

$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
